# "media dia 11 sept 2021"
#
# Sheet1 "ARQUITECTO": the voucher amount is bumped from $100,000 to
# $150,000, so both the numeric amount cell and its spelled-out text cell
# (amount-in-words) are updated, and this becomes the active sheet/selection
# (it was previously "VALES DE INSENTIVOS" that was active).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "ARQUITECTO        "

# Amount figure: 100,000 -> 150,000
$ws1.Range("D1").Value = 150000

# Amount in words, kept consistent with the new figure
$ws1.Range("A2").Value = "CIENTO CINCUENTA    MIL   PESOS 00/100 M.N."

# Make "ARQUITECTO" the active sheet/tab (previously "VALES DE INSENTIVOS"
# was the selected tab) and move its selection onto the amount block.
$ws1.Activate()
$ws1.Range("A4:D5").Select() | Out-Null
